$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (A7) -> set its value (B7) to "true"
$ws.Range("B7").Value = "true"

# "Date" row (A8) -> update value (B8) to new timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
